$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.088.29"
$ws.Range("E2").Value = "  -1.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.395.52"
$ws.Range("E3").Value = "  -3.95%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.65"
$ws.Range("E5").Value = "  -3.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.36"
$ws.Range("E6").Value = "  -5.87%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.395.59"
$ws.Range("E8").Value = "  -3.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -1.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.08"
$ws.Range("E10").Value = "  -9.72%  "

$ws.Range("E11").Value = "  -10.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.371"
$ws.Range("E12").Value = "  -8.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.973.29"
$ws.Range("E13").Value = "  -3.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000177"
$ws.Range("E14").Value = "  -11.01%  "

$ws.Range("E15").Value = "  -1.88%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "26.02"
$ws.Range("E16").Value = "  -8.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "65.028.20"
$ws.Range("E17").Value = "  -1.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.406.07"
$ws.Range("E18").Value = "  -3.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.48"
$ws.Range("E19").Value = "  -14.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.81"
$ws.Range("E20").Value = "  -6.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.46"
$ws.Range("E21").Value = "  -6.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "378.75"
$ws.Range("E22").Value = "  -9.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.548"
$ws.Range("E23").Value = "  -8.57%  "

$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.80"
$ws.Range("E25").Value = "  -7.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.531.60"
$ws.Range("E26").Value = "  -3.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000103"
$ws.Range("E27").Value = "  -10.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  -10.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").Value = "  -9.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.98"
$ws.Range("E31").Value = "  -10.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.407.51"
$ws.Range("E32").Value = "  -3.57%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("E34").Value = "  -7.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "22.67"
$ws.Range("E35").Value = "  -6.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.83"
$ws.Range("E36").Value = "  -3.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.63"
$ws.Range("E37").Value = "  -12.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.14"
$ws.Range("E38").Value = "  -11.51%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.45"
$ws.Range("E39").Value = "  -7.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.64"
$ws.Range("E40").Value = "  -12.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0748"
$ws.Range("E41").Value = "  -8.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.806"
$ws.Range("E42").Value = "  -6.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.31"
$ws.Range("E43").Value = "  -4.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.30"
$ws.Range("E45").Value = "  -15.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.58"
$ws.Range("E46").Value = "  -10.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.09"
$ws.Range("E47").Value = "  +1.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.66"
$ws.Range("E48").Value = "  -6.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.42"
$ws.Range("E49").Value = "  -8.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.02"
$ws.Range("E50").Value = "  -14.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.151.59"
$ws.Range("E51").Value = "  -8.56%  "
